# Update the Q0..Q9 naive QoQ error matrix (rows 2-20) with the refreshed
# ifoCAST-matched error series, and trim each row to its new (one-shorter)
# diagonal length now that the evaluation includes one more quarter.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -2.520798295803975
$ws.Range("C2").Value = -1.209498772153247
$ws.Range("D2").Value = 0.6650949949593515
$ws.Range("E2").Value = -0.05858795499538469
$ws.Range("F2").Value = 0.09291507488738329
$ws.Range("G2").Value = -0.1622689761646768
$ws.Range("H2").Value = 0.6645657133825503
$ws.Range("I2").Value = 0.5658615947808475
$ws.Range("J2").Value = 0.2128828273432847
$ws.Range("K2").Value = 0.3606851936599405

# Row 3
$ws.Range("B3").Value = -0.801266318531477
$ws.Range("C3").Value = 0.9167180388457457
$ws.Range("D3").Value = -0.14049400881852
$ws.Range("E3").Value = 0.1680071242686689
$ws.Range("F3").Value = -0.056361815650953
$ws.Range("G3").Value = 0.7229823595412057
$ws.Range("H3").Value = 0.6323002786109391
$ws.Range("I3").Value = 0.2875647859541212
$ws.Range("J3").Value = 0.4309984855609192
$ws.Range("K3").Value = 0.5879651212282146

# Row 4
$ws.Range("B4").Value = 0.8112188269969305
$ws.Range("C4").Value = -0.2665992927838778
$ws.Range("D4").Value = 0.3362551929019006
$ws.Range("E4").Value = 0.004302276952929285
$ws.Range("F4").Value = 0.7341151989158385
$ws.Range("G4").Value = 0.6951630299378696
$ws.Range("H4").Value = 0.3449163137478644
$ws.Range("I4").Value = 0.4755036625652983
$ws.Range("J4").Value = 0.6390389922546331
$ws.Range("K4").Value = 0.1406021036577823

# Row 5
$ws.Range("B5").Value = -0.8710399126316165
$ws.Range("C5").Value = 0.2720518155472072
$ws.Range("D5").Value = 0.2958495850101461
$ws.Range("E5").Value = 0.6964083228299049
$ws.Range("F5").Value = 0.6865790098644085
$ws.Range("G5").Value = 0.4286573913452327
$ws.Range("H5").Value = 0.5103308074933433
$ws.Range("I5").Value = 0.6656033753069301
$ws.Range("J5").Value = 0.1861009444043812
$ws.Range("K5").Value = 0.4585159253925983

# Row 6
$ws.Range("B6").Value = 0.1844109867082959
$ws.Range("C6").Value = 0.2434303829804775
$ws.Range("D6").Value = 0.7550547198149107
$ws.Range("E6").Value = 0.6915081552015099
$ws.Range("F6").Value = 0.4175641082963374
$ws.Range("G6").Value = 0.5230975067184298
$ws.Range("H6").Value = 0.6743329343542892
$ws.Range("I6").Value = 0.188602888553233
$ws.Range("J6").Value = 0.4647918628642795
$ws.Range("K6").Value = 0.2901314376394514

# Row 7
$ws.Range("B7").Value = 0.2658583626134861
$ws.Range("C7").Value = 0.7576582593773593
$ws.Range("D7").Value = 0.6694404695474946
$ws.Range("E7").Value = 0.4115821279919825
$ws.Range("F7").Value = 0.5177542565360187
$ws.Range("G7").Value = 0.663964400983923
$ws.Range("H7").Value = 0.1801255664296092
$ws.Range("I7").Value = 0.4570195545574687
$ws.Range("J7").Value = 0.2815074641064828
$ws.Range("K7").Value = 0.589805676679554

# Row 8
$ws.Range("B8").Value = 0.7790828308345794
$ws.Range("C8").Value = 0.7896183183467069
$ws.Range("D8").Value = 0.3223877835781994
$ws.Range("E8").Value = 0.4880272526874408
$ws.Range("F8").Value = 0.6725905210138331
$ws.Range("G8").Value = 0.154185393004174
$ws.Range("H8").Value = 0.4342409915292412
$ws.Range("I8").Value = 0.267889428957865
$ws.Range("J8").Value = 0.5712848794943781
$ws.Range("K8").Value = -0.04166000770222517

# Row 9
$ws.Range("B9").Value = 0.7472976860263175
$ws.Range("C9").Value = 0.2931833097224661
$ws.Range("D9").Value = 0.4850066528825197
$ws.Range("E9").Value = 0.6537292374023527
$ws.Range("F9").Value = 0.1338542619389697
$ws.Range("G9").Value = 0.4197108360949934
$ws.Range("H9").Value = 0.251248383465597
$ws.Range("I9").Value = 0.553701034065628
$ws.Range("J9").Value = -0.05813389123226798
$ws.Range("K9").Value = 0.6310880987550094

# Row 10
$ws.Range("B10").Value = 0.6333810503403763
$ws.Range("C10").Value = 0.5603777146664568
$ws.Range("D10").Value = 0.4626359510893536
$ws.Range("E10").Value = 0.1585481745257029
$ws.Range("F10").Value = 0.4290902396512893
$ws.Range("G10").Value = 0.196437584854206
$ws.Range("H10").Value = 0.5348069431769582
$ws.Range("I10").Value = -0.07299789204731372
$ws.Range("J10").Value = 0.6023763881256388
$ws.Range("K10").Value = 0.3304428394308724

# Row 11
$ws.Range("B11").Value = 1.011072586656707
$ws.Range("C11").Value = 0.509746653301712
$ws.Range("D11").Value = -0.08173505626510533
$ws.Range("E11").Value = 0.4650353207171521
$ws.Range("F11").Value = 0.1943868274146131
$ws.Range("G11").Value = 0.458238703740928
$ws.Range("H11").Value = -0.1003972864225971
$ws.Range("I11").Value = 0.5760873379279133
$ws.Range("J11").Value = 0.2870455487175398
$ws.Range("K11").ClearContents()

# Row 12
$ws.Range("B12").Value = 0.8220722402252505
$ws.Range("C12").Value = 0.051553683470419
$ws.Range("D12").Value = 0.2859568989009605
$ws.Range("E12").Value = 0.2225810661851803
$ws.Range("F12").Value = 0.4942640149333215
$ws.Range("G12").Value = -0.1379388566268107
$ws.Range("H12").Value = 0.5705449988405521
$ws.Range("I12").Value = 0.2912208776562884
$ws.Range("J12").ClearContents()

# Row 13
$ws.Range("B13").Value = 0.2871441745782602
$ws.Range("C13").Value = 0.3706356397752701
$ws.Range("D13").Value = 0.07698069955146319
$ws.Range("E13").Value = 0.5042177577925642
$ws.Range("F13").Value = -0.1225736869272658
$ws.Range("G13").Value = 0.5316966638831291
$ws.Range("H13").Value = 0.2753750686291025
$ws.Range("I13").ClearContents()

# Row 14
$ws.Range("B14").Value = 0.6816598262566529
$ws.Range("C14").Value = 0.1940694739626584
$ws.Range("D14").Value = 0.3417536163429973
$ws.Range("E14").Value = -0.0936147364620642
$ws.Range("F14").Value = 0.5674395363380327
$ws.Range("G14").Value = 0.2440474222454754
$ws.Range("H14").ClearContents()

# Row 15
$ws.Range("B15").Value = 0.4408946513667728
$ws.Range("C15").Value = 0.3594094838808868
$ws.Range("D15").Value = -0.1882369755730587
$ws.Range("E15").Value = 0.5996569245865127
$ws.Range("F15").Value = 0.2564355480731927
$ws.Range("G15").ClearContents()

# Row 16
$ws.Range("B16").Value = 0.5990858432970987
$ws.Range("C16").Value = -0.1032127321038452
$ws.Range("D16").Value = 0.4831723462284986
$ws.Range("E16").Value = 0.2715408197250452
$ws.Range("F16").ClearContents()

# Row 17
$ws.Range("B17").Value = 0.0616473449302421
$ws.Range("C17").Value = 0.4967096184764148
$ws.Range("D17").Value = 0.2085679007350822
$ws.Range("E17").ClearContents()

# Row 18
$ws.Range("B18").Value = 0.7505586603418228
$ws.Range("C18").Value = 0.3078859509171186
$ws.Range("D18").ClearContents()

# Row 19
$ws.Range("B19").Value = 0.3519456421565676
$ws.Range("C19").ClearContents()

# Row 20
$ws.Range("B20").ClearContents()
